$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 352
$ws.Range("I2").Value = 984
$ws.Range("J2").Value = 4180
$ws.Range("K2").Value = 16
$ws.Range("L2").Value = 1039
$ws.Range("M2").Value = 66
$ws.Range("N2").Value = 737
$ws.Range("P2").Value = 27
$ws.Range("Q2").Value = 14
$ws.Range("R2").Value = 49
$ws.Range("S2").Value = 465
$ws.Range("T2").Value = 738
$ws.Range("U2").Value = 62
$ws.Range("V2").Value = 6462
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 6386
$ws.Range("Y2").Value = 14
$ws.Range("Z2").Value = 90
$ws.Range("AA2").Value = 42
